$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 was missing its "great" flag in column G - fill it in (same style
# and "no" value used elsewhere in that column).
$ws.Range("G25").Value = "no"

# Copy the formatting of an existing "appid/keyword" row (row 2) down onto
# the two new rows (26-27) so the new rows pick up the same cell styles
# (bold/centered appid+recovery columns etc.) as the rest of the table.
$ws.Range("A2:G2").Copy()
$ws.Range("A26:G27").PasteSpecial(-4122)

# Row 26 - new review
$ws.Range("A26").Value = "com.hamxa.shaynachim"
$ws.Range("B26").Value = "bitcoin"
$ws.Range("C26").Value = "sofershani9@gmail.com"
$ws.Range("D26").Value = "rotemzinger3@gmail.com"
$ws.Range("E26").Value = "27/5/2019 15:59"
$ws.Range("F26").Value = "I hope you are ready for this great guide"
$ws.Range("G26").Value = "yes"

# Row 27 - new review
$ws.Range("A27").Value = "com.hamxa.shaynachim"
$ws.Range("B27").Value = "bitcoin"
$ws.Range("C27").Value = "emmakrigel63@gmail.com "
$ws.Range("D27").Value = "sofershani9@gmail.com"
$ws.Range("E27").Value = "27/5/2019 15:59"
$ws.Range("F27").Value = "excellent"
$ws.Range("G27").Value = "yes"

# New hyperlink on the recovery-email cell of row 26 (same address as the
# hyperlink already used on D22).
$ws.Hyperlinks.Add($ws.Range("D26"), "mailto:rotemzinger3@gmail.com", $null, $null, "rotemzinger3@gmail.com")

# Adding the hyperlink re-styles the cell with the default blue/underlined
# "Hyperlink" look; re-apply the plain table formatting (taken from D22,
# which already carries a hyperlink but keeps the normal cell style) so the
# cell format matches the rest of the "recovery" column.
$ws.Range("D22").Copy()
$ws.Range("D26").PasteSpecial(-4122)

# Move the selection to the first empty row below the newly added data,
# matching where the author's cursor ended up after the paste.
$ws.Range("A28").Select() | Out-Null
